# Update the K column (G) values for peralta_wandy.xlsx
# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 0
    3 = 2
    4 = 0
    5 = 1
    6 = 3
    8 = 1
    9 = 0
    12 = 3
    13 = 1
    14 = 2
    15 = 2
    16 = 3
    17 = 0
    19 = 2
    20 = 0
    21 = 0
    23 = 1
    24 = 1
    25 = 1
    26 = 0
    27 = 2
    28 = 0
    29 = 1
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    37 = 0
    38 = 0
    39 = 3
    40 = 0
    41 = 1
    42 = 1
    43 = 0
    44 = 2
    46 = 0
    47 = 1
    48 = 2
    49 = 0
    50 = 2
    51 = 0
    52 = 0
    53 = 1
    54 = 0
    55 = 0
    56 = 0
    58 = 2
    59 = 2
    60 = 0
    61 = 2
    62 = 1
    63 = 2
    64 = 3
    65 = 0
    66 = 1
    67 = 0
    68 = 2
    69 = 1
    70 = 2
    71 = 3
    72 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
